# Scheduled-runner data refresh: updates the current market-board prices
# (and the derived Leve-profit figures in columns H-N) for each crafting
# class's table -- ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 55: H55, I55, K55, M55
$ws.Range("H55").Value2 = 2456.7144
$ws.Range("I55").Value2 = 436.25
$ws.Range("K55").Value2 = 436.25
$ws.Range("M55").Value2 = -222.25
# Row 96: H96, I96, K96, M96
$ws.Range("H96").Value2 = 654.6
$ws.Range("I96").Value2 = 654.6
$ws.Range("K96").Value2 = 1963.8
$ws.Range("M96").Value2 = -590.8000000000002
# Row 98: H98, I98, J98, K98, L98, M98, N98
$ws.Range("H98").Value2 = 3056.9302
$ws.Range("I98").Value2 = 3036.35
$ws.Range("J98").Value2 = 3331.3333
$ws.Range("K98").Value2 = 3036.35
$ws.Range("L98").Value2 = 3331.3333
$ws.Range("M98").Value2 = -1538.35
$ws.Range("N98").Value2 = -6327.3333
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value2 = 320.2857
$ws.Range("I107").Value2 = 293.0909
$ws.Range("J107").Value2 = 420
$ws.Range("K107").Value2 = 293.0909
$ws.Range("L107").Value2 = 420
$ws.Range("M107").Value2 = 1626.9091
$ws.Range("N107").Value2 = -4260
# Row 113: H113, I113, K113, M113
$ws.Range("H113").Value2 = 90911410
$ws.Range("I113").Value2 = 125002420
$ws.Range("K113").Value2 = 125002420
$ws.Range("M113").Value2 = -124999166
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value2 = 3056.9302
$ws.Range("I122").Value2 = 3036.35
$ws.Range("J122").Value2 = 3331.3333
$ws.Range("K122").Value2 = 9109.049999999999
$ws.Range("L122").Value2 = 9993.999899999999
$ws.Range("M122").Value2 = -6659.049999999999
$ws.Range("N122").Value2 = -14893.9999
# Row 127: H127, I127, K127, M127
$ws.Range("H127").Value2 = 647
$ws.Range("I127").Value2 = 647
$ws.Range("K127").Value2 = 1941
$ws.Range("M127").Value2 = 3019
# Row 137: H137, I137, K137, M137
$ws.Range("H137").Value2 = 2133.8572
$ws.Range("I137").Value2 = 2109.25
$ws.Range("K137").Value2 = 6327.75
$ws.Range("M137").Value2 = -3777.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45: H45, I45, K45, M45
$ws.Range("H45").Value2 = 55025.156
$ws.Range("I45").Value2 = 75963.25999999999
$ws.Range("K45").Value2 = 75963.25999999999
$ws.Range("M45").Value2 = -75586.25999999999
# Row 74: H74, I74, K74, M74
$ws.Range("H74").Value2 = 9144
$ws.Range("I74").Value2 = 7470.857
$ws.Range("K74").Value2 = 7470.857
$ws.Range("M74").Value2 = -6596.857
# Row 77: H77, I77, K77, M77
$ws.Range("H77").Value2 = 9144
$ws.Range("I77").Value2 = 7470.857
$ws.Range("K77").Value2 = 37354.285
$ws.Range("M77").Value2 = -32986.285
# Row 93: H93, J93, L93; clear N93
$ws.Range("H93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("N93").ClearContents()
# Row 97: H97, I97, K97, M97
$ws.Range("H97").Value2 = 513.1818
$ws.Range("I97").Value2 = 349.44446
$ws.Range("K97").Value2 = 349.44446
$ws.Range("M97").Value2 = 146.55554
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value2 = 3880.8572
$ws.Range("I122").Value2 = 3630.182
$ws.Range("K122").Value2 = 10890.546
$ws.Range("M122").Value2 = -8440.545999999998
# Row 131: H131, J131, L131; clear N131
$ws.Range("H131").Value2 = 0
$ws.Range("J131").Value2 = 0
$ws.Range("L131").Value2 = 0
$ws.Range("N131").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105: H105, I105, J105, K105, L105, M105, N105
$ws.Range("H105").Value2 = 5554.85
$ws.Range("I105").Value2 = 5207.2144
$ws.Range("J105").Value2 = 6366
$ws.Range("K105").Value2 = 5207.2144
$ws.Range("L105").Value2 = 6366
$ws.Range("M105").Value2 = -3460.2144
$ws.Range("N105").Value2 = -9860
# Row 107: H107, I107, K107, M107
$ws.Range("H107").Value2 = 4550.75
$ws.Range("I107").Value2 = 3486.5715
$ws.Range("K107").Value2 = 3486.5715
$ws.Range("M107").Value2 = -1566.5715
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value2 = 5702.625
$ws.Range("I134").Value2 = 5842.1304
$ws.Range("J134").Value2 = 2494
$ws.Range("K134").Value2 = 17526.3912
$ws.Range("L134").Value2 = 7482
$ws.Range("M134").Value2 = -14991.3912
$ws.Range("N134").Value2 = -12552

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value2 = 5275.1665
$ws.Range("I58").Value2 = 3750
$ws.Range("J58").Value2 = 6037.75
$ws.Range("K58").Value2 = 3750
$ws.Range("L58").Value2 = 6037.75
$ws.Range("M58").Value2 = -3547
$ws.Range("N58").Value2 = -6443.75
# Row 74: H74, J74, L74, N74
$ws.Range("H74").Value2 = 41909.6
$ws.Range("J74").Value2 = 41909.6
$ws.Range("L74").Value2 = 41909.6
$ws.Range("N74").Value2 = -43657.6
# Row 77: H77, J77, L77, N77
$ws.Range("H77").Value2 = 41909.6
$ws.Range("J77").Value2 = 41909.6
$ws.Range("L77").Value2 = 125728.8
$ws.Range("N77").Value2 = -134464.8
# Row 94: H94, J94, L94, N94
$ws.Range("H94").Value2 = 3536.3076
$ws.Range("J94").Value2 = 3758.6
$ws.Range("L94").Value2 = 3758.6
$ws.Range("N94").Value2 = -4660.6
# Row 99: H99, I99, K99, M99
$ws.Range("H99").Value2 = 7944.8
$ws.Range("I99").Value2 = 6974.909
$ws.Range("K99").Value2 = 6974.909
$ws.Range("M99").Value2 = -5476.909
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value2 = 1616.5217
$ws.Range("I107").Value2 = 556.9167
$ws.Range("J107").Value2 = 2772.4546
$ws.Range("K107").Value2 = 556.9167
$ws.Range("L107").Value2 = 2772.4546
$ws.Range("M107").Value2 = 1363.0833
$ws.Range("N107").Value2 = -6612.4546
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value2 = 7944.8
$ws.Range("I126").Value2 = 6974.909
$ws.Range("K126").Value2 = 20924.727
$ws.Range("M126").Value2 = -18454.727
# Row 130: H130, J130, L130, N130
$ws.Range("H130").Value2 = 94896.5
$ws.Range("J130").Value2 = 94896.5
$ws.Range("L130").Value2 = 94896.5
$ws.Range("N130").Value2 = -104936.5
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value2 = 3590.25
$ws.Range("I132").Value2 = 2184.3333
$ws.Range("J132").Value2 = 4996.1665
$ws.Range("K132").Value2 = 6552.999899999999
$ws.Range("L132").Value2 = 14988.4995
$ws.Range("M132").Value2 = -4022.999899999999
$ws.Range("N132").Value2 = -20048.4995
# Row 134: H134, J134, L134, N134
$ws.Range("H134").Value2 = 6657.636
$ws.Range("J134").Value2 = 13124.75
$ws.Range("L134").Value2 = 39374.25
$ws.Range("N134").Value2 = -44444.25
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value2 = 5275.1665
$ws.Range("I136").Value2 = 3750
$ws.Range("J136").Value2 = 6037.75
$ws.Range("K136").Value2 = 11250
$ws.Range("L136").Value2 = 18113.25
$ws.Range("M136").Value2 = -8700
$ws.Range("N136").Value2 = -23213.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122: H122, J122, L122, N122
$ws.Range("H122").Value2 = 2925
$ws.Range("J122").Value2 = 3600
$ws.Range("L122").Value2 = 32400
$ws.Range("N122").Value2 = -37300
# Row 134: H134, I134, K134, M134
$ws.Range("H134").Value2 = 1147
$ws.Range("I134").Value2 = 1147
$ws.Range("K134").Value2 = 3441
$ws.Range("M134").Value2 = 1629
# Row 140: H140, I140, K140, M140
$ws.Range("H140").Value2 = 3461.7778
$ws.Range("I140").Value2 = 3523.5557
$ws.Range("K140").Value2 = 10570.6671
$ws.Range("M140").Value2 = -5390.667099999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70, I70, J70, K70, L70, M70, N70
$ws.Range("H70").Value2 = 17107.46
$ws.Range("I70").Value2 = 16399.666
$ws.Range("J70").Value2 = 17714.143
$ws.Range("K70").Value2 = 16399.666
$ws.Range("L70").Value2 = 17714.143
$ws.Range("M70").Value2 = -16129.666
$ws.Range("N70").Value2 = -18254.143
# Row 73: H73, I73, J73, K73, L73, M73, N73
$ws.Range("H73").Value2 = 17107.46
$ws.Range("I73").Value2 = 16399.666
$ws.Range("J73").Value2 = 17714.143
$ws.Range("K73").Value2 = 16399.666
$ws.Range("L73").Value2 = 17714.143
$ws.Range("M73").Value2 = -15463.666
$ws.Range("N73").Value2 = -19586.143
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value2 = 1983.1666
$ws.Range("I122").Value2 = 1824.75
$ws.Range("J122").Value2 = 2300
$ws.Range("K122").Value2 = 5474.25
$ws.Range("L122").Value2 = 6900
$ws.Range("M122").Value2 = -3024.25
$ws.Range("N122").Value2 = -11800

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16, J16, L16, N16
$ws.Range("H16").Value2 = 2385.158
$ws.Range("J16").Value2 = 932
$ws.Range("L16").Value2 = 932
$ws.Range("N16").Value2 = -1272
# Row 20: H20, J20, L20, N20
$ws.Range("H20").Value2 = 15998.75
$ws.Range("J20").Value2 = 17331.666
$ws.Range("L20").Value2 = 17331.666
$ws.Range("N20").Value2 = -17783.666
# Row 40: H40, I40, J40, K40, L40, M40, N40
$ws.Range("H40").Value2 = 2926.2
$ws.Range("I40").Value2 = 3140.3333
$ws.Range("J40").Value2 = 999
$ws.Range("K40").Value2 = 3140.3333
$ws.Range("L40").Value2 = 999
$ws.Range("M40").Value2 = -3004.3333
$ws.Range("N40").Value2 = -1271
# Row 68: H68, J68, L68, N68
$ws.Range("H68").Value2 = 9221.333000000001
$ws.Range("J68").Value2 = 8669
$ws.Range("L68").Value2 = 8669
$ws.Range("N68").Value2 = -10167
# Row 71: H71, J71, L71, N71
$ws.Range("H71").Value2 = 9221.333000000001
$ws.Range("J71").Value2 = 8669
$ws.Range("L71").Value2 = 43345
$ws.Range("N71").Value2 = -50833
# Row 82: H82, I82, K82, M82
$ws.Range("H82").Value2 = 4080.92
$ws.Range("I82").Value2 = 2535.6
$ws.Range("K82").Value2 = 2535.6
$ws.Range("M82").Value2 = -2174.6
# Row 85: H85, I85, K85, M85
$ws.Range("H85").Value2 = 4080.92
$ws.Range("I85").Value2 = 2535.6
$ws.Range("K85").Value2 = 2535.6
$ws.Range("M85").Value2 = -1287.6
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value2 = 4334.231
$ws.Range("I136").Value2 = 3919
$ws.Range("K136").Value2 = 11757
$ws.Range("M136").Value2 = -9207

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 15: H15, J15, L15, N15
$ws.Range("H15").Value2 = 7583.3335
$ws.Range("J15").Value2 = 8750
$ws.Range("L15").Value2 = 8750
$ws.Range("N15").Value2 = -9326
# Row 113: H113, I113, K113, M113
$ws.Range("H113").Value2 = 540.7273
$ws.Range("I113").Value2 = 593
$ws.Range("K113").Value2 = 1779
$ws.Range("M113").Value2 = 391
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value2 = 5177.5557
$ws.Range("I126").Value2 = 5177.5557
$ws.Range("K126").Value2 = 15532.6671
$ws.Range("M126").Value2 = -13062.6671
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value2 = 5116.2285
$ws.Range("I136").Value2 = 4123.6924
$ws.Range("J136").Value2 = 7983.5557
$ws.Range("K136").Value2 = 12371.0772
$ws.Range("L136").Value2 = 23950.6671
$ws.Range("M136").Value2 = -9821.0772
$ws.Range("N136").Value2 = -29050.6671

